$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("заказы")

$ws.Range("F1").Value = "Дата создания"
$ws.Range("F2").Value = "2024 01 08 04:43:42"
$ws.Range("F3").Value = "2024 01 08 04:44:01"
$ws.Range("F4").Value = "2024 01 08 16:54:34"
